# Actualización desde MV -datos-
# Append the September/October 2021 daily auction rows to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115
$ws.Range("A115").Value = "20-09-2021"
$ws.Range("B115").Value = 800000
$ws.Range("C115").Value = 895000
$ws.Range("D115").Value = 400000
$ws.Range("E115").Value = 345000
$ws.Range("F115").Value = 55000
$ws.Range("G115").Value = 1.59

# Row 116
$ws.Range("A116").Value = "21-09-2021"
$ws.Range("B116").Value = 800000
$ws.Range("C116").Value = 720000
$ws.Range("D116").Value = 400000
$ws.Range("E116").Value = 250000
$ws.Range("F116").Value = 150000
$ws.Range("G116").Value = 1.53

# Row 117
$ws.Range("A117").Value = "22-09-2021"
$ws.Range("B117").Value = 1000000
$ws.Range("D117").Value = 0

# Row 118
$ws.Range("A118").Value = "23-09-2021"
$ws.Range("B118").Value = 800000
$ws.Range("D118").Value = 0

# Row 119
$ws.Range("A119").Value = "24-09-2021"
$ws.Range("B119").Value = 800000
$ws.Range("D119").Value = 0

# Row 120
$ws.Range("A120").Value = "27-09-2021"
$ws.Range("B120").Value = 300000
$ws.Range("D120").Value = 0

# Row 121
$ws.Range("A121").Value = "28-09-2021"
$ws.Range("B121").Value = 300000
$ws.Range("D121").Value = 0

# Row 122
$ws.Range("A122").Value = "29-09-2021"
$ws.Range("B122").Value = 300000
$ws.Range("C122").Value = 450000
$ws.Range("D122").Value = 180000
$ws.Range("E122").Value = 180000
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 1.9

# Row 123
$ws.Range("A123").Value = "30-09-2021"
$ws.Range("B123").Value = 300000
$ws.Range("D123").Value = 0

# Row 124
# "01-10-2021" looks like a valid date (day=01 <= 12 could be a month), so a
# plain .Value assignment gets auto-converted to a date serial by Excel's
# smart-parsing. Force it to stay literal text (matches the source data,
# which stores it as a shared string) by temporarily marking the cell as
# Text before assigning, then resetting the cell style back to Normal so no
# visible formatting is left behind.
$ws.Range("A124").NumberFormat = "@"
$ws.Range("A124").Value = "01-10-2021"
$ws.Range("A124").Style = "Normal"
$ws.Range("B124").Value = 300000
$ws.Range("D124").Value = 0
